$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain plain text so purely-numeric-looking values
# (e.g. "545.33") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.559.99'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '3.083.14'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '545.33'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").Value = '140.10'
$ws.Range("E6").Value = '  +1.67%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '3.076.92'
$ws.Range("E8").Value = '  -0.81%  '
$ws.Range("D9").Value = '0.504'
$ws.Range("E9").Value = '  +1.28%  '
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").Value = '6.39'
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("E12").Value = '  -2.79%  '
$ws.Range("D13").Value = '35.07'
$ws.Range("E13").Value = '  -1.62%  '
$ws.Range("D14").Value = '0.0000225'
$ws.Range("E14").Value = '  +2.92%  '
$ws.Range("D15").Value = '3.585.07'
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("D16").Value = '63.559.61'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").Value = '3.084.14'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("E19").Value = '  -1.38%  '
$ws.Range("D20").Value = '475.37'
$ws.Range("E20").Value = '  -3.46%  '
$ws.Range("D21").Value = '13.48'
$ws.Range("E21").Value = '  -1.39%  '
$ws.Range("D22").Value = '0.702'
$ws.Range("E22").Value = '  -2.71%  '
$ws.Range("E23").Value = '  -2.19%  '
$ws.Range("D24").Value = '78.85'
$ws.Range("E24").Value = '  -0.54%  '
$ws.Range("E25").Value = '  -1.19%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").Value = '7.98'
$ws.Range("E28").Value = '  -6.05%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").Value = '26.25'
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("E31").Value = '  -3.81%  '
$ws.Range("E32").Value = '  +1.87%  '
$ws.Range("D33").Value = '57.91'
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("E34").Value = '  -7.77%  '
$ws.Range("D35").Value = '5.43'
$ws.Range("E35").Value = '  +5.25%  '
$ws.Range("D36").Value = '493.58'
$ws.Range("E36").Value = '  -4.99%  '
$ws.Range("D37").Value = '6.02'
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").Value = '3.247.78'
$ws.Range("E38").Value = '  +2.85%  '
$ws.Range("D39").Value = '0.0404'
$ws.Range("E39").Value = '  -0.50%  '
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("D41").Value = '0.118'
$ws.Range("E41").Value = '  -1.72%  '
$ws.Range("D42").Value = '8.14'
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("D43").Value = '2.64'
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("D44").Value = '0.255'
$ws.Range("E44").Value = '  -1.98%  '
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = '124.95'
$ws.Range("E46").Value = '  +3.03%  '
$ws.Range("D47").Value = '25.50'
$ws.Range("E47").Value = '  +1.42%  '
$ws.Range("E48").Value = '  -1.73%  '
$ws.Range("E49").Value = '  +3.88%  '
$ws.Range("E50").Value = '  +0.81%  '
$ws.Range("D51").Value = '2.31'
$ws.Range("E51").Value = '  -0.34%  '

# Restore the default (no explicit format override) style on column D
# so the saved file matches the original formatting of these cells.
$ws.Range("D2:D51").Style = "Normal"
